# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.513.66"
$ws.Range("E2").Value = "  +2.13%  "

$ws.Range("D3").Value = "3.561.20"
$ws.Range("E3").Value = "  +1.12%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'598.14"
$ws.Range("E5").Value = "  +1.89%  "

$ws.Range("D6").Value = "'172.22"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("D7").Value = "3.555.70"
$ws.Range("E7").Value = "  +1.10%  "

$ws.Range("D8").Value = "'0.614"
$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("E10").Value = "  +4.51%  "

$ws.Range("D11").Value = "'7.41"
$ws.Range("E11").Value = "  +9.74%  "

$ws.Range("D12").Value = "'0.586"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").Value = "'46.31"
$ws.Range("E13").Value = "  -2.41%  "

$ws.Range("D15").Value = "4.137.06"
$ws.Range("E15").Value = "  +1.25%  "

$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").Value = "'610.74"

$ws.Range("D18").Value = "3.563.45"
$ws.Range("E18").Value = "  +1.35%  "

$ws.Range("D19").Value = "70.551.82"
$ws.Range("E19").Value = "  +2.18%  "

$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").Value = "'17.33"
$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("E22").Value = "  -0.34%  "

$ws.Range("D23").Value = "'9.23"
$ws.Range("E23").Value = "  -16.76%  "

$ws.Range("D24").Value = "'15.74"
$ws.Range("E24").Value = "  -0.01%  "

$ws.Range("D25").Value = "'96.68"
$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("E26").Value = "  -2.51%  "

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").Value = "'2.60"
$ws.Range("E28").Value = "  -0.18%  "

$ws.Range("D29").Value = "'33.79"
$ws.Range("E29").Value = "  +3.73%  "

$ws.Range("D30").Value = "'9.04"
$ws.Range("E30").Value = "  -1.64%  "

$ws.Range("E31").Value = "  -2.84%  "

$ws.Range("E32").Value = "  -2.45%  "

$ws.Range("D33").Value = "'661.05"
$ws.Range("E33").Value = "  +8.07%  "

$ws.Range("D34").Value = "'7.09"
$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").Value = "'1.29"
$ws.Range("E35").Value = "  -1.49%  "

$ws.Range("D36").Value = "'3.60"
$ws.Range("E36").Value = "  +4.18%  "

$ws.Range("D37").Value = "'0.100"
$ws.Range("E37").Value = "  -1.57%  "

$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'57.35"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0473"
$ws.Range("E40").Value = "  +6.89%  "

$ws.Range("E41").Value = "  +0.08%  "

$ws.Range("E42").Value = "  +4.34%  "

$ws.Range("D43").Value = "3.378.98"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("D44").Value = "'0.319"
$ws.Range("E44").Value = "  -1.75%  "

$ws.Range("D45").Value = "0.0₃0703"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").Value = "'32.75"
$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").Value = "'2.95"
$ws.Range("E47").Value = "  +7.51%  "

$ws.Range("D48").Value = "'2.62"
$ws.Range("E48").Value = "  +4.18%  "

$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("D50").Value = "'132.37"
$ws.Range("E50").Value = "  -0.97%  "

$ws.Range("E51").Value = "  -0.09%  "
